$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell format (style s="1") from H1 into the new header cells
# so the new columns match the existing header formatting exactly.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows for the two new columns (I: rows 2-11, J: rows 2-11)
$iValues = @(4, 5, 1, 1, 1, 2, 1, 1, 3, 1)
$jValues = @(7, 6, 3, 5, 5, 5, 4, 2, 4, 2)

for ($r = 0; $r -lt 10; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
